# Reverse the order of comma-separated names in column G ("Recorded By")
# for every data row on the active sheet. Cells with a single value are
# left unchanged (reversing a one-element list is a no-op).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($null -ne $val -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ", "
        $reversed = $parts[($parts.Count - 1)..0]
        $cell.Value = [string]::Join(", ", $reversed)
    }
}
